$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chapter_4_Table_S4.11")

# --- Column A combination labels: underscore-joined gene symbols -> comma-separated,
# with updated gene names (MTM1->SMG8, ELF2->ELP2 etc reflecting corrected annotation) ---
$ws.Range("A4").Value = "SATB1, ALAS1"
$ws.Range("A5").Value = "SATB1, ALAS1, NR4A2"
$ws.Range("A6").Value = "SATB1, ALAS1, NR4A2, SMG8"
$ws.Range("A7").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2"
$ws.Range("A8").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5"
$ws.Range("A9").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D"
$ws.Range("A10").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3"
$ws.Range("A11").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3, MYB"
$ws.Range("A12").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3, MYB, BCDIN3D"
$ws.Range("A13").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3, MYB, BCDIN3D, TEX2"
$ws.Range("A14").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3, MYB, BCDIN3D, TEX2, MAD2L2"
$ws.Range("A15").Value = "SATB1, ALAS1, NR4A2, SMG8, ELP2, DCAF5, ZNF75D, PRICKLE3, MYB, BCDIN3D, TEX2, MAD2L2, KATNBL1"

# --- Column I (Group) label: "Pass 1" -> "13-gene set" for all data rows ---
$ws.Range("I3").Value = "13-gene set"
$ws.Range("I4").Value = "13-gene set"
$ws.Range("I5").Value = "13-gene set"
$ws.Range("I6").Value = "13-gene set"
$ws.Range("I7").Value = "13-gene set"
$ws.Range("I8").Value = "13-gene set"
$ws.Range("I9").Value = "13-gene set"
$ws.Range("I10").Value = "13-gene set"
$ws.Range("I11").Value = "13-gene set"
$ws.Range("I12").Value = "13-gene set"
$ws.Range("I13").Value = "13-gene set"
$ws.Range("I14").Value = "13-gene set"
$ws.Range("I15").Value = "13-gene set"

# --- Updated cross-validation fold scores / averages (re-run results) ---
$ws.Range("D3").Value = 0.61328125
$ws.Range("G3").Value = 0.79608916567639199
$ws.Range("B6").Value = 0.96598639455782298
$ws.Range("C6").Value = 0.97213622291021695
$ws.Range("D6").Value = 0.91015625
$ws.Range("E6").Value = 0.97499999999999998
$ws.Range("F6").Value = 0.94230769230769196
$ws.Range("G6").Value = 0.95311731195514604
$ws.Range("B7").Value = 0.97959183673469397
$ws.Range("C7").Value = 0.98142414860681104
$ws.Range("D7").Value = 0.9453125
$ws.Range("F7").Value = 0.93589743589743601
$ws.Range("G7").Value = 0.95927851758112104
$ws.Range("B8").Value = 0.97959183673469397
$ws.Range("C8").Value = 0.98142414860681104
$ws.Range("D8").Value = 0.94921875
$ws.Range("E8").Value = 0.98333333333333295
$ws.Range("F8").Value = 0.96153846153846201
$ws.Range("G8").Value = 0.97102130604266002
$ws.Range("B9").Value = 0.98299319727891199
$ws.Range("C9").Value = 0.98761609907120695
$ws.Range("D9").Value = 0.95703125
$ws.Range("E9").Value = 0.97916666666666696
$ws.Range("F9").Value = 0.96153846153846201
$ws.Range("G9").Value = 0.97366913491104901
$ws.Range("B10").Value = 0.98639455782312901
$ws.Range("C10").Value = 0.97523219814241502
$ws.Range("D10").Value = 0.96484375
$ws.Range("E10").Value = 0.97916666666666696
$ws.Range("F10").Value = 0.99358974358974395
$ws.Range("G10").Value = 0.97984538324439097
$ws.Range("B11").Value = 0.98299319727891199
$ws.Range("C11").Value = 0.98761609907120695
$ws.Range("D11").Value = 0.98828125
$ws.Range("E11").Value = 0.98750000000000004
$ws.Range("F11").Value = 0.98076923076923095
$ws.Range("G11").Value = 0.98543195542386997
$ws.Range("C12").Value = 0.97832817337461297
$ws.Range("D12").Value = 1
$ws.Range("G12").Value = 0.99059760746403802
$ws.Range("B13").Value = 0.98979591836734704
$ws.Range("C13").Value = 0.99380804953560398
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.98333333333333295
$ws.Range("G13").Value = 0.99338746024725699
$ws.Range("B14").Value = 0.98979591836734704
$ws.Range("C14").Value = 0.98761609907120695
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.99583333333333302
$ws.Range("G14").Value = 0.99464907015437798
$ws.Range("B15").Value = 0.98979591836734704
$ws.Range("C15").Value = 0.99380804953560398
$ws.Range("D15").Value = 0.99609375
$ws.Range("G15").Value = 0.99593954358059

# --- Row heights: data rows grew from 20.1pt to 24.95pt (to match header row) ---
$ws.Rows(2).RowHeight = 24.95
$ws.Rows(3).RowHeight = 24.95
$ws.Rows(4).RowHeight = 24.95
$ws.Rows(5).RowHeight = 24.95
$ws.Rows(6).RowHeight = 24.95
$ws.Rows(7).RowHeight = 24.95
$ws.Rows(8).RowHeight = 24.95
$ws.Rows(9).RowHeight = 24.95
$ws.Rows(10).RowHeight = 24.95
$ws.Rows(11).RowHeight = 24.95
$ws.Rows(12).RowHeight = 24.95
$ws.Rows(13).RowHeight = 24.95
$ws.Rows(14).RowHeight = 24.95
$ws.Rows(15).RowHeight = 24.95

# --- Column widths: widened to fit the new, longer text ---
$ws.Columns(1).ColumnWidth = 94.85546875
$ws.Columns(9).ColumnWidth = 10.85546875
